$wb = $excel.ActiveWorkbook

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1489.6818
$ws.Range("I100").Value = 589
$ws.Range("J100").Value = 2790.6667
$ws.Range("K100").Value = 589
$ws.Range("L100").Value = 2790.6667
$ws.Range("M100").Value = -48
$ws.Range("N100").Value = -3872.6667

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3917.4285
$ws.Range("I125").Value = 4686.4
$ws.Range("J125").Value = 1995
$ws.Range("K125").Value = 42177.6
$ws.Range("L125").Value = 17955
$ws.Range("M125").Value = -39717.6
$ws.Range("N125").Value = -22875

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 20001362
$ws.Range("I61").Value = 20001362
$ws.Range("K61").Value = 20001362
$ws.Range("M61").Value = -20001150

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 16131439
$ws.Range("I74").Value = 23810718
$ws.Range("J74").Value = 4951.4
$ws.Range("K74").Value = 23810718
$ws.Range("L74").Value = 4951.4
$ws.Range("M74").Value = -23809844
$ws.Range("N74").Value = -6699.4

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 16131439
$ws.Range("I77").Value = 23810718
$ws.Range("J77").Value = 4951.4
$ws.Range("K77").Value = 119053590
$ws.Range("L77").Value = 24757
$ws.Range("M77").Value = -119049222
$ws.Range("N77").Value = -33493

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8066427.5
$ws.Range("I132").Value = 13159228
$ws.Range("J132").Value = 2826
$ws.Range("K132").Value = 39477684
$ws.Range("L132").Value = 8478
$ws.Range("M132").Value = -39475154
$ws.Range("N132").Value = -13538

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 20001362
$ws.Range("I136").Value = 20001362
$ws.Range("K136").Value = 60004086
$ws.Range("M136").Value = -60001536

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1617.7142
$ws.Range("I107").Value = 1456
$ws.Range("K107").Value = 1456
$ws.Range("M107").Value = 464

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2878.6274
$ws.Range("I134").Value = 2195.325
$ws.Range("J134").Value = 5363.364
$ws.Range("K134").Value = 6585.974999999999
$ws.Range("L134").Value = 16090.092
$ws.Range("M134").Value = -4050.974999999999
$ws.Range("N134").Value = -21160.092

# CRP row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 40775.8
$ws.Range("J140").Value = 40775.8
$ws.Range("L140").Value = 40775.8
$ws.Range("N140").Value = -51135.8

# CUL row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 9813.538
$ws.Range("J97").Value = 2437.3
$ws.Range("L97").Value = 7311.900000000001
$ws.Range("N97").Value = -8303.900000000001

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 17991
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -22891

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 813.9091
$ws.Range("I131").Value = 356.66666
$ws.Range("J131").Value = 859.63336
$ws.Range("K131").Value = 1069.99998
$ws.Range("L131").Value = 2578.90008
$ws.Range("M131").Value = 3970.00002
$ws.Range("N131").Value = -12658.90008

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 805.05554
$ws.Range("I132").Value = 533.9167
$ws.Range("J132").Value = 1347.3334
$ws.Range("K132").Value = 4805.2503
$ws.Range("L132").Value = 12126.0006
$ws.Range("M132").Value = -2275.2503
$ws.Range("N132").Value = -17186.0006

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4945.7646
$ws.Range("I133").Value = 3719.7778
$ws.Range("J133").Value = 6325
$ws.Range("K133").Value = 11159.3334
$ws.Range("L133").Value = 18975
$ws.Range("M133").Value = -6099.3334
$ws.Range("N133").Value = -29095

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 69666.664

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 69666.664

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4072.9473
$ws.Range("I102").Value = 5292
$ws.Range("J102").Value = 1431.6666
$ws.Range("K102").Value = 5292
$ws.Range("L102").Value = 1431.6666
$ws.Range("M102").Value = -3670
$ws.Range("N102").Value = -4675.6666

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2736.2173
$ws.Range("I132").Value = 1797.6129
$ws.Range("J132").Value = 4676
$ws.Range("K132").Value = 5392.8387
$ws.Range("L132").Value = 14028
$ws.Range("M132").Value = -2862.8387
$ws.Range("N132").Value = -19088

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3830.7173
$ws.Range("I7").Value = 3634.6956
$ws.Range("K7").Value = 3634.6956
$ws.Range("M7").Value = -3522.6956

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1142
$ws.Range("I22").Value = 766.6667
$ws.Range("J22").Value = 1282.75
$ws.Range("K22").Value = 766.6667
$ws.Range("L22").Value = 1282.75
$ws.Range("M22").Value = -471.6667
$ws.Range("N22").Value = -1872.75

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1142
$ws.Range("I27").Value = 766.6667
$ws.Range("J27").Value = 1282.75
$ws.Range("K27").Value = 766.6667
$ws.Range("L27").Value = 1282.75
$ws.Range("M27").Value = -659.6667
$ws.Range("N27").Value = -1496.75

# LTW row 30
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1600
$ws.Range("I30").Value = 1600
$ws.Range("K30").Value = 1600
$ws.Range("M30").Value = -1492

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1378.5834
$ws.Range("I61").Value = 1258.5883
$ws.Range("J61").Value = 1670
$ws.Range("K61").Value = 1258.5883
$ws.Range("L61").Value = 1670
$ws.Range("M61").Value = -1056.5883
$ws.Range("N61").Value = -2074

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1378.5834
$ws.Range("I113").Value = 1258.5883
$ws.Range("J113").Value = 1670
$ws.Range("K113").Value = 1258.5883
$ws.Range("L113").Value = 1670
$ws.Range("M113").Value = 911.4117000000001
$ws.Range("N113").Value = -6010

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3830.7173
$ws.Range("I126").Value = 3634.6956
$ws.Range("K126").Value = 10904.0868
$ws.Range("M126").Value = -8434.086800000001

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7515.5117
$ws.Range("I132").Value = 5370.737
$ws.Range("J132").Value = 9213.458000000001
$ws.Range("K132").Value = 16112.211
$ws.Range("L132").Value = 27640.374
$ws.Range("M132").Value = -13582.211
$ws.Range("N132").Value = -32700.374

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1116.0555
$ws.Range("I132").Value = 595.0769
$ws.Range("J132").Value = 2470.6
$ws.Range("K132").Value = 1785.2307
$ws.Range("L132").Value = 7411.799999999999
$ws.Range("M132").Value = 744.7692999999999
$ws.Range("N132").Value = -12471.8

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1575.381
$ws.Range("I136").Value = 1204.15
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 3612.45
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -1062.45
$ws.Range("N136").Value = -32100
